$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.347.51'
$ws.Range('E2').Value = '  -1.81%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.887.01'
$ws.Range('E3').Value = '  -1.94%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9980'
$ws.Range('E4').Value = '  -0.41%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '236.82'
$ws.Range('E5').Value = '  -1.62%  '

$ws.Range('B6').Value = 'USDC'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9985'
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4839'
$ws.Range('E7').Value = '  -1.50%  '

$ws.Range('B8').Value = 'Cardano'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2890'
$ws.Range('E8').Value = '  -3.03%  '

$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06609'
$ws.Range('E9').Value = '  -2.51%  '

$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.869.71'
$ws.Range('E10').Value = '  -2.77%  '

$ws.Range('B11').Value = 'Solana'
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '16.94'
$ws.Range('E11').Value = '  -1.01%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07335'
$ws.Range('E12').Value = '  +0.46%  '

$ws.Range('B13').Value = 'Polkadot'
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.131'
$ws.Range('E13').Value = '  -0.95%  '

$ws.Range('B14').Value = 'Litecoin'
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '87.69'
$ws.Range('E14').Value = '  -2.30%  '

$ws.Range('B15').Value = 'Polygon'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6613'
$ws.Range('E15').Value = '  -1.81%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.295.67'
$ws.Range('E16').Value = '  -1.87%  '

$ws.Range('B17').Value = 'Avalanche'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.39'
$ws.Range('E17').Value = '  -1.30%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007771'
$ws.Range('E18').Value = '  -2.90%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.9986'
$ws.Range('E19').Value = '  -0.31%  '

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.403'
$ws.Range('E20').Value = '  +4.13%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '2.123.36'
$ws.Range('E21').Value = '  -1.88%  '

$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9984'
$ws.Range('E22').Value = '  -0.39%  '

$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '195.15'
$ws.Range('E23').Value = '  -5.47%  '

$ws.Range('B24').Value = 'Chainlink'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.183'
$ws.Range('E24').Value = '  -2.10%  '

$ws.Range('B25').Value = 'Cosmos'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.297'
$ws.Range('E25').Value = '  -4.22%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.45'
$ws.Range('E26').Value = '  +3.31%  '

$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.15'
$ws.Range('E27').Value = '  -5.19%  '

$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.932'
$ws.Range('E28').Value = '  -3.11%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.305'
$ws.Range('E30').Value = '  -1.46%  '

$ws.Range('B31').Value = 'Stellar'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09153'
$ws.Range('E31').Value = '  -0.47%  '

$ws.Range('B32').Value = 'Filecoin'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.032'
$ws.Range('E32').Value = '  -1.46%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05069'
$ws.Range('E33').Value = '  -2.62%  '

$ws.Range('B34').Value = 'ARBITRUM'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.149'
$ws.Range('E34').Value = '  +1.92%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7234'
$ws.Range('E35').Value = '  -4.34%  '

$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.693'
$ws.Range('E36').Value = '  -1.34%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01789'
$ws.Range('E37').Value = '  -3.96%  '

$ws.Range('B38').Value = 'MXToken'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.638'
$ws.Range('E38').Value = '  -3.71%  '

$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9190'
$ws.Range('E39').Value = '  -1.02%  '

$ws.Range('B40').Value = 'RenderToken'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.057'
$ws.Range('E40').Value = '  -1.97%  '

$ws.Range('B41').Value = 'Quant'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '106.12'
$ws.Range('E41').Value = '  -1.98%  '

$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4312'
$ws.Range('E42').Value = '  -4.85%  '

$ws.Range('B43').Value = 'FraxShare'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.816'
$ws.Range('E43').Value = '  -2.27%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.001'
$ws.Range('E44').Value = '  -1.03%  '

$ws.Range('B45').Value = 'Aptos'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '7.483'
$ws.Range('E45').Value = '  -2.93%  '

$ws.Range('B46').Value = 'Algorand'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.1315'
$ws.Range('E46').Value = '  -6.04%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.563'
$ws.Range('E47').Value = '  +8.12%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '64.99'
$ws.Range('E48').Value = '  -8.45%  '

$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.854'
$ws.Range('E49').Value = '  -2.75%  '

$ws.Range('B50').Value = 'Elrond'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '33.93'
$ws.Range('E50').Value = '  -4.62%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.05743'
$ws.Range('E51').Value = '  -3.58%  '
